$d = $word.ActiveDocument

# We need to add 5 new paragraphs at the very top of the document body:
#   1. "Introduction"                (bold, underlined heading)
#   2. body paragraph describing the document
#   3. empty paragraph (spacer)
#   4. "Content" (bold, underlined heading) followed by a trailing empty run
#   5. empty paragraph (spacer)
#
# Strategy: insert 6 blank paragraphs before the current first paragraph,
# fill paragraphs 1, 2 and 4 with text/formatting, then merge the 6th
# (still-empty) paragraph into paragraph 4 by deleting the paragraph mark
# between them -- this produces the trailing, unformatted empty run that
# follows the bold "Content" run in paragraph 4, while paragraphs 3 and 5
# remain as untouched empty paragraphs.

$insertionPoint = $d.Range(0, 0)
for ($i = 0; $i -lt 6; $i++) {
    $insertionPoint.InsertParagraphBefore()
}

# Paragraph 1: "Introduction" heading
$para1 = $d.Paragraphs(1).Range
$para1.Text = "Introduction"
$para1.Bold = 1
$para1.Font.Underline = 1

# Paragraph 2: introductory body text
$apostrophe = [char]0x2019
$para2 = $d.Paragraphs(2).Range
$para2.Text = "This document is an exploration of how OpenAI logs its requests with language models they host. The reason this was investigated is that we weren" + $apostrophe + "t aware of batch requests at first, and thought there could be an issue in the prompt. "

# Paragraph 3 is left as an empty spacer paragraph.

# Paragraph 4: "Content" heading, with a trailing empty (unformatted) run
$para4 = $d.Paragraphs(4).Range
$para4.Text = "Content"
$para4.Bold = 1
$para4.Font.Underline = 1

# Merge paragraph 5 (still empty) into paragraph 4 so its empty run becomes
# a second, non-bold/non-underlined run trailing the "Content" run.
$paraMark = $d.Range($d.Paragraphs(4).Range.End - 1, $d.Paragraphs(4).Range.End)
$paraMark.Delete()

# Paragraph 5 (was paragraph 6) is left as an empty spacer paragraph.

Write-Output "Paragraphs after edit: $($d.Paragraphs.Count)"
